$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.941.46'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '3.699.46'
$ws.Range('E3').Value = '  +2.80%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '240.17'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').Value = '1.89'
$ws.Range('E6').Value = '  +8.68%  '
$ws.Range('D7').Value = '655.00'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -2.12%  '
$ws.Range('E9').Value = '  +3.36%  '
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').Value = '3.695.96'
$ws.Range('E11').Value = '  +2.68%  '
$ws.Range('D12').Value = '45.45'
$ws.Range('E12').Value = '  +2.18%  '
$ws.Range('D13').Value = '0.206'
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('D14').Value = '6.90'
$ws.Range('E14').Value = '  +6.43%  '
$ws.Range('D15').Value = '4.385.62'
$ws.Range('E15').Value = '  +2.83%  '
$ws.Range('E16').Value = '  +2.33%  '
$ws.Range('D17').Value = '96.698.17'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').Value = '9.10'
$ws.Range('E18').Value = '  +4.86%  '
$ws.Range('D19').Value = '3.700.42'
$ws.Range('E19').Value = '  +2.53%  '
$ws.Range('D20').Value = '19.34'
$ws.Range('E20').Value = '  +6.37%  '
$ws.Range('D21').Value = '12.90'
$ws.Range('E21').Value = '  +2.37%  '
$ws.Range('D22').Value = '0.528'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').Value = '526.78'
$ws.Range('E23').Value = '  +1.53%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').Value = '7.14'
$ws.Range('E25').Value = '  +2.76%  '
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('D27').Value = '101.99'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').Value = '13.44'
$ws.Range('E28').Value = '  +2.22%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '3.893.87'
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.169'
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '12.61'
$ws.Range('E31').Value = '  +4.32%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '3.05'
$ws.Range('E32').Value = '  +2.17%  '
$ws.Range('B33').Value = 'Dai'
$ws.Range('C33').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '1.87'
$ws.Range('E34').Value = '  +14.05%  '
$ws.Range('B35').Value = 'Cronos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D35').Value = '0.186'
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = '32.78'
$ws.Range('E36').Value = '  +2.37%  '
$ws.Range('B37').Value = 'Binance-PegBSC-USD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').Value = '0.611'
$ws.Range('E38').Value = '  +6.49%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '656.87'
$ws.Range('E39').Value = '  +5.83%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').Value = '9.06'
$ws.Range('E40').Value = '  +3.48%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '7.08'
$ws.Range('E41').Value = '  +17.18%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.164'
$ws.Range('E42').Value = '  +5.57%  '
$ws.Range('B43').Value = 'ImmutableX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D43').Value = '2.02'
$ws.Range('E43').Value = '  +3.53%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '0.969'
$ws.Range('E44').Value = '  +3.88%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '38.83'
$ws.Range('E45').Value = '  +17.93%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.457'
$ws.Range('E47').Value = '  +5.68%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0461'
$ws.Range('E48').Value = '  +4.05%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '2.33'
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '8.82'
$ws.Range('E50').Value = '  +2.86%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '23.63'
$ws.Range('E51').Value = '  -0.10%  '
